# "Added tests to complianceReportsPageTest4":
# insert a new TC154/A row right under the header on the CRPT sheet
# (pushing the existing rows down by one), and leave the workbook
# focused on the CRPT sheet/cell A26 instead of CRPTEthane.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRPT")

$ws.Rows.Item(2).Insert() | Out-Null
$ws.Rows.Item(2).ClearFormats() | Out-Null
$ws.Cells.Item(2, 1).Value = "TC154"
$ws.Cells.Item(2, 2).Value = "A"
$ws.Cells.Item(2, 3).Value = $null

$ws.Select() | Out-Null
$ws.Range("A26").Select() | Out-Null
